$wb = $excel.ActiveWorkbook

# --- About sheet: update the "last updated" date from 3/15/2024 to 3/28/2024 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- RAF-capacity sheet: hydrogen combustion turbine / combined cycle RAF raised to 1 ---
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# --- Update the active sheet / selection to RAF-capacity, matching the saved view state ---
$wsCapacity.Activate()
$wsCapacity.Range("B25").Select()
$excel.ActiveWindow.Zoom = 80
